# Generate Report for Handback
#
# The handback status report gained a second processed file
# (3c891a4d-1a77-43ce-b189-283a1c6f85c5.md) alongside the original
# (which itself got re-run, landing a new handoff uuid:
# 0f4edcdc-4655-4d1c-8b4a-aca8da87b601.md). Update row 2 on every
# sheet with the fresh identifiers/timestamps and append a new row 3
# for the newly tracked file.

$wb = $excel.ActiveWorkbook

function Set-Text($range, [string]$text) {
    # Force text storage (avoid Excel auto-typing "True"/"False" as
    # booleans or otherwise reinterpreting the literal).
    $range.Value = "'" + $text
}

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ovTable = $ov.ListObjects.Item(1)

# Row 2: bd0234ce... -> 0f4edcdc... (re-handed-back) + refreshed date
Set-Text $ov.Range("A2") "0f4edcdc-4655-4d1c-8b4a-aca8da87b601.md"
$ov.Range("B2").Value = "e2e\0f4edcdc-4655-4d1c-8b4a-aca8da87b601.md"
Set-Text $ov.Range("G2") "2016-08-26 13:02:11"

# New row 3 for 3c891a4d...
$ovTable.ListRows.Add() | Out-Null
Set-Text $ov.Range("A3") "3c891a4d-1a77-43ce-b189-283a1c6f85c5.md"
$ov.Range("B3").Value = "e2e\3c891a4d-1a77-43ce-b189-283a1c6f85c5.md"
Set-Text $ov.Range("C3") ".md"
Set-Text $ov.Range("E3") "Handed back: in sync with en-US"
Set-Text $ov.Range("F3") "Handed back: in sync with en-US"
Set-Text $ov.Range("G3") "2016-08-26 13:02:11"

$ov.Range("B3").Font.Underline = $true
$ov.Range("B3").Font.Color = 15570276
$ov.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77c8b3d4c33f1fbae2da4857990b128aecccd04e/e2e/0f4edcdc-4655-4d1c-8b4a-aca8da87b601.md", "", "", "e2e\0f4edcdc-4655-4d1c-8b4a-aca8da87b601.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77c8b3d4c33f1fbae2da4857990b128aecccd04e/e2e/3c891a4d-1a77-43ce-b189-283a1c6f85c5.md", "", "", "e2e\3c891a4d-1a77-43ce-b189-283a1c6f85c5.md") | Out-Null

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zhTable = $zh.ListObjects.Item(1)

Set-Text $zh.Range("A2") "0f4edcdc-4655-4d1c-8b4a-aca8da87b601.md"
Set-Text $zh.Range("G2") "0f4edcdc-4655-4d1c-8b4a-aca8da87b601.b98abda519c198c0511e765eb7c698bf1f7292a9.zh-cn.xlf"
Set-Text $zh.Range("H2") "2016-08-26 13:01:58"
Set-Text $zh.Range("I2") "0f4edcdc-4655-4d1c-8b4a-aca8da87b601.md"
Set-Text $zh.Range("J2") "0f4edcdc-4655-4d1c-8b4a-aca8da87b601.b98abda519c198c0511e765eb7c698bf1f7292a9.zh-cn.xlf"
Set-Text $zh.Range("K2") "2016-08-26 13:02:32"

$zhTable.ListRows.Add() | Out-Null
Set-Text $zh.Range("A3") "3c891a4d-1a77-43ce-b189-283a1c6f85c5.md"
Set-Text $zh.Range("B3") ".md"
Set-Text $zh.Range("C3") "Handed back: in sync with en-US"
Set-Text $zh.Range("D3") "e2e"
Set-Text $zh.Range("E3") "ht"
Set-Text $zh.Range("F3") "True"
Set-Text $zh.Range("G3") "3c891a4d-1a77-43ce-b189-283a1c6f85c5.3a5e08c4d9328c2a861dce6b4e5f08ee29823cbd.zh-cn.xlf"
Set-Text $zh.Range("H3") "2016-08-26 13:01:58"
Set-Text $zh.Range("I3") "3c891a4d-1a77-43ce-b189-283a1c6f85c5.md"
Set-Text $zh.Range("J3") "3c891a4d-1a77-43ce-b189-283a1c6f85c5.3a5e08c4d9328c2a861dce6b4e5f08ee29823cbd.zh-cn.xlf"
Set-Text $zh.Range("K3") "2016-08-26 13:02:32"
Set-Text $zh.Range("M3") "True"
Set-Text $zh.Range("O3") "False"

$zh.Range("A3").Font.Underline = $true
$zh.Range("A3").Font.Color = 15570276
$zh.Range("I3").Font.Underline = $true
$zh.Range("I3").Font.Color = 15570276
$zh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77c8b3d4c33f1fbae2da4857990b128aecccd04e/e2e/3c891a4d-1a77-43ce-b189-283a1c6f85c5.md", "", "", "3c891a4d-1a77-43ce-b189-283a1c6f85c5.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/667bd41c4dc96ad63b1f94eb6b2164329700c559/e2e/3c891a4d-1a77-43ce-b189-283a1c6f85c5.md", "", "", "3c891a4d-1a77-43ce-b189-283a1c6f85c5.md") | Out-Null

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$deTable = $de.ListObjects.Item(1)

Set-Text $de.Range("A2") "0f4edcdc-4655-4d1c-8b4a-aca8da87b601.md"
Set-Text $de.Range("G2") "0f4edcdc-4655-4d1c-8b4a-aca8da87b601.b98abda519c198c0511e765eb7c698bf1f7292a9.de-de.xlf"
Set-Text $de.Range("I2") "0f4edcdc-4655-4d1c-8b4a-aca8da87b601.md"
Set-Text $de.Range("J2") "0f4edcdc-4655-4d1c-8b4a-aca8da87b601.b98abda519c198c0511e765eb7c698bf1f7292a9.de-de.xlf"
Set-Text $de.Range("K2") "2016-08-26 13:02:39"

$deTable.ListRows.Add() | Out-Null
Set-Text $de.Range("A3") "3c891a4d-1a77-43ce-b189-283a1c6f85c5.md"
Set-Text $de.Range("B3") ".md"
Set-Text $de.Range("C3") "Handed back: in sync with en-US"
Set-Text $de.Range("D3") "e2e"
Set-Text $de.Range("E3") "ht"
Set-Text $de.Range("F3") "True"
Set-Text $de.Range("G3") "3c891a4d-1a77-43ce-b189-283a1c6f85c5.3a5e08c4d9328c2a861dce6b4e5f08ee29823cbd.de-de.xlf"
Set-Text $de.Range("H3") "2016-08-26 13:02:11"
Set-Text $de.Range("I3") "3c891a4d-1a77-43ce-b189-283a1c6f85c5.md"
Set-Text $de.Range("J3") "3c891a4d-1a77-43ce-b189-283a1c6f85c5.3a5e08c4d9328c2a861dce6b4e5f08ee29823cbd.de-de.xlf"
Set-Text $de.Range("K3") "2016-08-26 13:02:39"
Set-Text $de.Range("M3") "True"
Set-Text $de.Range("O3") "False"

$de.Range("A3").Font.Underline = $true
$de.Range("A3").Font.Color = 15570276
$de.Range("I3").Font.Underline = $true
$de.Range("I3").Font.Color = 15570276
$de.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77c8b3d4c33f1fbae2da4857990b128aecccd04e/e2e/3c891a4d-1a77-43ce-b189-283a1c6f85c5.md", "", "", "3c891a4d-1a77-43ce-b189-283a1c6f85c5.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/ad74ca710a5a446445ac72106cff467cc641290c/e2e/3c891a4d-1a77-43ce-b189-283a1c6f85c5.md", "", "", "3c891a4d-1a77-43ce-b189-283a1c6f85c5.md") | Out-Null

Write-Output "handback status report updated"
